# S28/G01: Backtesting foundation (page + runs + EOD candles)
# Appends sprint S28 task rows (281-301) to the sprint tasks tracker sheet,
# covering groups G01-G06 for the new Backtesting epic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TaskCell {
    param($Row, $Col, $Text)
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.Value() = $Text
    $cell.WrapText() = $false
    $cell.VerticalAlignment() = -4107
}

# Row 281
Set-TaskCell 281 1 "S28"
Set-TaskCell 281 2 "G01"
Set-TaskCell 281 3 "Backtesting foundation (page + runs + EOD data)"
Set-TaskCell 281 4 "S28_G01_TF001"
Set-TaskCell 281 5 "Frontend: Add Backtesting page in sidebar (below Alerts) with 3 tabs (Signal/Portfolio/Execution) and a split layout (Inputs vs Results)."
Set-TaskCell 281 7 "implemented"
Set-TaskCell 281 8 "Backtesting page added (/backtesting) with 3 tabs + help dialog."

# Row 282
Set-TaskCell 282 1 "S28"
Set-TaskCell 282 2 "G01"
Set-TaskCell 282 3 "Backtesting foundation (page + runs + EOD data)"
Set-TaskCell 282 4 "S28_G01_TB001"
Set-TaskCell 282 5 "Backend: Add DB schema for backtest runs (config snapshot + timestamps + status) and API endpoints to create/list/get runs."
Set-TaskCell 282 7 "implemented"
Set-TaskCell 282 8 "Backtest run table + CRUD endpoints added."

# Row 283
Set-TaskCell 283 1 "S28"
Set-TaskCell 283 2 "G01"
Set-TaskCell 283 3 "Backtesting foundation (page + runs + EOD data)"
Set-TaskCell 283 4 "S28_G01_TB002"
Set-TaskCell 283 5 "Backend: Implement EOD candle loader for a selected universe (Holdings/Group/Both) with strict “as-of” semantics (no lookahead) and caching."
Set-TaskCell 283 7 "implemented"
Set-TaskCell 283 8 "EOD candle loader implemented (aligned close matrix) + endpoint."

# Row 284
Set-TaskCell 284 1 "S28"
Set-TaskCell 284 2 "G01"
Set-TaskCell 284 3 "Backtesting foundation (page + runs + EOD data)"
Set-TaskCell 284 4 "S28_G01_TF002"
Set-TaskCell 284 5 "Frontend: Add “Runs” panel (recent runs/history), ability to rerun with the same config, and basic compare (A vs B) scaffolding."
Set-TaskCell 284 7 "implemented"
Set-TaskCell 284 8 "Runs panel added (history + rerun scaffolding)."

# Row 285
Set-TaskCell 285 1 "S28"
Set-TaskCell 285 2 "G01"
Set-TaskCell 285 3 "Backtesting foundation (page + runs + EOD data)"
Set-TaskCell 285 4 "S28_G01_TT001"
Set-TaskCell 285 5 "Tests: Add backend tests for backtest run CRUD + candle loader invariants (no-lookahead, deterministic results for fixed inputs)."
Set-TaskCell 285 7 "implemented"
Set-TaskCell 285 8 "Backend tests added for runs + candle loader."

# Row 286
Set-TaskCell 286 1 "S28"
Set-TaskCell 286 2 "G02"
Set-TaskCell 286 3 "Signal backtest (EOD)"
Set-TaskCell 286 4 "S28_G02_TB001"
Set-TaskCell 286 5 "Backend: Implement Signal backtest endpoint (DSL condition and/or ranking) producing hit-rate + forward-return distributions (1D/5D/20D)."
Set-TaskCell 286 7 "planned"

# Row 287
Set-TaskCell 287 1 "S28"
Set-TaskCell 287 2 "G02"
Set-TaskCell 287 3 "Signal backtest (EOD)"
Set-TaskCell 287 4 "S28_G02_TF001"
Set-TaskCell 287 5 "Frontend: Signal backtest tab UI (universe, date range, DSL/ranking input, forward windows) + results (summary + distributions)."
Set-TaskCell 287 7 "planned"

# Row 288
Set-TaskCell 288 1 "S28"
Set-TaskCell 288 2 "G02"
Set-TaskCell 288 3 "Signal backtest (EOD)"
Set-TaskCell 288 4 "S28_G02_TT001"
Set-TaskCell 288 5 "Tests: Validate signal backtest metrics on small synthetic datasets and guard against lookahead bias."
Set-TaskCell 288 7 "planned"

# Row 289
Set-TaskCell 289 1 "S28"
Set-TaskCell 289 2 "G03"
Set-TaskCell 289 3 "Portfolio backtest v1 (target weights)"
Set-TaskCell 289 4 "S28_G03_TB001"
Set-TaskCell 289 5 "Backend: Implement portfolio simulator (EOD) with rebalance cadence, budget/max-trades/min-trade constraints, and simple costs/slippage model."
Set-TaskCell 289 7 "planned"

# Row 290
Set-TaskCell 290 1 "S28"
Set-TaskCell 290 2 "G03"
Set-TaskCell 290 3 "Portfolio backtest v1 (target weights)"
Set-TaskCell 290 4 "S28_G03_TB002"
Set-TaskCell 290 5 "Backend: Implement Target-weights portfolio backtest (use portfolio group target weights) returning equity curve, drawdowns, turnover, and rebalance actions."
Set-TaskCell 290 7 "planned"

# Row 291
Set-TaskCell 291 1 "S28"
Set-TaskCell 291 2 "G03"
Set-TaskCell 291 3 "Portfolio backtest v1 (target weights)"
Set-TaskCell 291 4 "S28_G03_TF001"
Set-TaskCell 291 5 "Frontend: Portfolio backtest tab (Target weights mode) + results views (equity/drawdown charts, turnover, action list, contributors)."
Set-TaskCell 291 7 "planned"

# Row 292
Set-TaskCell 292 1 "S28"
Set-TaskCell 292 2 "G03"
Set-TaskCell 292 3 "Portfolio backtest v1 (target weights)"
Set-TaskCell 292 4 "S28_G03_TT001"
Set-TaskCell 292 5 "Tests: Portfolio simulator invariants (cash/position accounting, constraints respected, deterministic outputs)."
Set-TaskCell 292 7 "planned"

# Row 293
Set-TaskCell 293 1 "S28"
Set-TaskCell 293 2 "G04"
Set-TaskCell 293 3 "Portfolio backtest v2 (rotation)"
Set-TaskCell 293 4 "S28_G04_TB001"
Set-TaskCell 293 5 "Backend: Implement rotation target derivation (Top-N selection + weighting) using DSL/strategy outputs and eligibility filters; plug into portfolio backtest."
Set-TaskCell 293 7 "planned"

# Row 294
Set-TaskCell 294 1 "S28"
Set-TaskCell 294 2 "G04"
Set-TaskCell 294 3 "Portfolio backtest v2 (rotation)"
Set-TaskCell 294 4 "S28_G04_TF001"
Set-TaskCell 294 5 "Frontend: Rotation mode UI (ranking source, Top-N, weighting, cadence) + compare vs Target-weights baseline."
Set-TaskCell 294 7 "planned"

# Row 295
Set-TaskCell 295 1 "S28"
Set-TaskCell 295 2 "G04"
Set-TaskCell 295 3 "Portfolio backtest v2 (rotation)"
Set-TaskCell 295 4 "S28_G04_TT001"
Set-TaskCell 295 5 "Tests: Rotation selection stability + turnover bounds + filter correctness."
Set-TaskCell 295 7 "planned"

# Row 296
Set-TaskCell 296 1 "S28"
Set-TaskCell 296 2 "G05"
Set-TaskCell 296 3 "Portfolio backtest v3 (risk parity)"
Set-TaskCell 296 4 "S28_G05_TB001"
Set-TaskCell 296 5 "Backend: Implement risk-parity target derivation for backtests (window/lookback + constraints) using covariance/risk metrics with EOD candles."
Set-TaskCell 296 7 "planned"

# Row 297
Set-TaskCell 297 1 "S28"
Set-TaskCell 297 2 "G05"
Set-TaskCell 297 3 "Portfolio backtest v3 (risk parity)"
Set-TaskCell 297 4 "S28_G05_TF001"
Set-TaskCell 297 5 "Frontend: Risk parity mode UI (window, constraints) + reporting focused on risk-adjusted outcomes and drawdowns."
Set-TaskCell 297 7 "planned"

# Row 298
Set-TaskCell 298 1 "S28"
Set-TaskCell 298 2 "G05"
Set-TaskCell 298 3 "Portfolio backtest v3 (risk parity)"
Set-TaskCell 298 4 "S28_G05_TT001"
Set-TaskCell 298 5 "Tests: Risk parity weight sanity checks and stability (constraints applied, non-negative weights, sum-to-1 within tolerance)."
Set-TaskCell 298 7 "planned"

# Row 299
Set-TaskCell 299 1 "S28"
Set-TaskCell 299 2 "G06"
Set-TaskCell 299 3 "Execution backtest (EOD friction model)"
Set-TaskCell 299 4 "S28_G06_TB001"
Set-TaskCell 299 5 "Backend: Add execution friction layer (fill timing: close vs next open, slippage bps, simple charges) to portfolio backtests."
Set-TaskCell 299 7 "planned"

# Row 300
Set-TaskCell 300 1 "S28"
Set-TaskCell 300 2 "G06"
Set-TaskCell 300 3 "Execution backtest (EOD friction model)"
Set-TaskCell 300 4 "S28_G06_TF001"
Set-TaskCell 300 5 "Frontend: Execution backtest tab (select base portfolio config + friction knobs) with “ideal vs realistic” comparison."
Set-TaskCell 300 7 "planned"

# Row 301
Set-TaskCell 301 1 "S28"
Set-TaskCell 301 2 "G06"
Set-TaskCell 301 3 "Execution backtest (EOD friction model)"
Set-TaskCell 301 4 "S28_G06_TT001"
Set-TaskCell 301 5 "Tests: Execution friction impacts (cost reduces returns; no negative-cash violations; deterministic)."
Set-TaskCell 301 7 "planned"

